$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I16:I40").Value = 3.447169811320755
